$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.820.00"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "1.640.82"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  -0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.76"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.506"
$ws.Range("E6").Value = "  +1.59%  "
$ws.Range("E7").Value = "  -0.64%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.252"
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0620"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.72"
$ws.Range("E10").Value = "  +4.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "1.870.26"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").Value = "1.650.72"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.527"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.32"
$ws.Range("E16").Value = "  +3.84%  "
$ws.Range("D17").Value = "26.841.86"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.43"
$ws.Range("E19").Value = "  +3.96%  "
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.37"
$ws.Range("E21").Value = "  +1.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.60"
$ws.Range("E22").Value = "  +7.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.42"
$ws.Range("E23").Value = "  +4.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.14"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.50"
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.35"
$ws.Range("E27").Value = "  +5.20%  "
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.82"
$ws.Range("E29").Value = "  +2.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0510"
$ws.Range("E30").Value = "  +2.25%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.35"
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.99"
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.54"
$ws.Range("E34").Value = "  +2.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.45"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").Value = "1.236.77"
$ws.Range("E36").Value = "  -1.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0173"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.538"
$ws.Range("E38").Value = "  +3.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.832"
$ws.Range("E39").Value = "  +4.29%  "
$ws.Range("E40").Value = "  -0.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.804"
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.36"
$ws.Range("E42").Value = "  +2.48%  "
$ws.Range("D43").Value = "1.783.43"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("E44").Value = "  -3.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "60.64"
$ws.Range("E45").Value = "  +1.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.57"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("D48").Value = "0.0₆0105"
$ws.Range("E48").Value = "  +10.25%  "
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0970"
$ws.Range("E50").Value = "  +2.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.52"
$ws.Range("E51").Value = "  +1.38%  "
